$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two parameter rows: nu_leakSR (row 43) and g_PMCA (row 44)
# Before: A43=nu_leakSR/B43=0.2 , A44=g_PMCA/B44=5.37
# After:  A43=g_PMCA/B43=5.37   , A44=nu_leakSR/B44=0.2
$ws.Range("A43").Value = "g_PMCA"
$ws.Range("B43").Value = 5.37
$ws.Range("A44").Value = "nu_leakSR"
$ws.Range("B44").Value = 0.2

# Add a new parameter row for g_leakNa
$ws.Range("A45").Value = "g_leakNa"
$ws.Range("B45").Value = 0.002

# Update the view/selection to match the saved state
$ws.Range("B46").Select()
